# TrialsSetup workbook refresh (2025-12-22 12:00)
#
# The "ALLEGRETTO-LTE (B7981028)" trial row in the query-table output
# (Sheet1!A11:B11) picked up a refreshed value from the underlying Power
# Query / SharePoint list: the trial name lost its stray leading space and
# its Progress figure moved from 0 to 62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "ALLEGRETTO-LTE (B7981028)"
$ws.Range("B11").Value = 62
